$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values actually vary row-to-row in the data block (rows 100-182).
$cols = @("D","J","K","L","M","O","P")

# 1) Snapshot current values for rows 100..182 before writing anything,
#    so the shift-down doesn't clobber data it still needs to read.
$vals = @{}
foreach ($col in $cols) {
    $vals[$col] = @{}
    for ($r = 100; $r -le 182; $r++) {
        $vals[$col][$r] = $ws.Range("$col$r").Value2
    }
}

# 2) Shift every snapshot row down by one: new row (r+1) gets old row r's data.
#    This pushes old row 182 into new row 183 and old row 100 into new row 101.
for ($r = 182; $r -ge 100; $r--) {
    foreach ($col in $cols) {
        $ws.Range("$col$($r+1)").Value2 = $vals[$col][$r]
    }
}

# 3) Fill the rest of new row 183 (the columns that never change) by copying
#    from row 182, which used to hold that data and is constant across rows anyway.
$ws.Range("A183").Value2 = $ws.Range("A182").Value2
$ws.Range("B183").Value2 = $ws.Range("B182").Value2
$ws.Range("C183").Value2 = $ws.Range("C182").Value2
$ws.Range("E183").Value2 = $ws.Range("E182").Value2
$ws.Range("F183").Value2 = $ws.Range("F182").Value2
$ws.Range("G183").Value2 = $ws.Range("G182").Value2
$ws.Range("H183").Value2 = $ws.Range("H182").Value2
$ws.Range("I183").Value2 = $ws.Range("I182").Value2
$ws.Range("N183").Value2 = $ws.Range("N182").Value2
$ws.Range("Q183").Value2 = $ws.Range("Q182").Value2
$ws.Range("R183").Value2 = $ws.Range("R182").Value2

# D is date-formatted (numFmt 165) in every data row; give the new row 183's
# D cell the same number format so it carries the style, not just the value.
$ws.Range("D183").NumberFormat = $ws.Range("D182").NumberFormat

# 4) New data point inserted at row 100 (pushing the old rows down, as above).
$ws.Range("D100").Value2 = 44447
$ws.Range("J100").Value2 = 800
$ws.Range("K100").Value2 = 4500
$ws.Range("L100").Value2 = 5000
$ws.Range("M100").Value2 = 4750
$ws.Range("P100").Value2 = 238
